$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content first
$ws.Range("A1:C20").ClearContents()

$ws.Range("A1").Value = "building"
$ws.Range("B1").Value = "units being built"
$ws.Range("C1").Value = "time left"
$ws.Range("A2").Value = "Nexus"
$ws.Range("A3").Value = "Gateway"
$ws.Range("B3").Formula = "'19"
$ws.Range("C3").Formula = "'5"
$ws.Range("A4").Value = "Gateway"
$ws.Range("A5").Value = "Forge"
$ws.Range("A6").Value = "Nexus"
$ws.Range("A7").Value = "Nexus"
$ws.Range("A8").Value = "Nexus"
$ws.Range("A9").Value = "Cibernetics Core"
$ws.Range("B9").Value = "43, 46"
$ws.Range("C9").Value = "85,  129"
$ws.Range("A10").Value = "Forge"
$ws.Range("A11").Value = "Forge"
$ws.Range("B11").Formula = "'37"
$ws.Range("C11").Formula = "'0"
$ws.Range("A12").Value = "Stargate"
$ws.Range("B12").Value = "20, 22"
$ws.Range("C12").Value = "0,  37"
$ws.Range("A13").Value = "Nexus"
$ws.Range("A14").Value = "Nexus"
$ws.Range("A15").Value = "Nexus"
$ws.Range("A16").Value = "Cibernetics Core"
$ws.Range("A17").Value = "Cibernetics Core"
$ws.Range("A18").Value = "Gateway"
$ws.Range("A19").Value = "Gateway"
$ws.Range("A20").Value = "Cibernetics Core"
$ws.Range("A21").Value = "Stargate"
$ws.Range("B21").Formula = "'20"
$ws.Range("C21").Formula = "'1"
$ws.Range("A22").Value = "Nexus"
$ws.Range("A23").Value = "Forge"

# Remove the quote-prefix styling applied for text-forced numeric cells
$ws.Range("B3").ClearFormats()
$ws.Range("C3").ClearFormats()
$ws.Range("B11").ClearFormats()
$ws.Range("C11").ClearFormats()
$ws.Range("B21").ClearFormats()
$ws.Range("C21").ClearFormats()
